$d = $word.ActiveDocument

# Locate the semicolon that terminates the "...n DEG 137/2017;" sentence and
# collapse to the insertion point right before it (i.e. right after "2017").
$anchor = $d.Content
$found = $anchor.Find.Execute("n" + [char]0x00B0 + " 137/2017;", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the target sentence."
}
$anchor.Collapse(0) | Out-Null
$anchor.MoveStart(1, -1) | Out-Null
$insertPos = $anchor.Start

# Insert the new middle clause as its own run, right before the old ";".
$newText = " e aggiornato con delibera n" + [char]0x00B0 + " 177/2025, Verb. 521"
$midRange = $d.Range($insertPos, $insertPos)
$midRange.InsertBefore($newText) | Out-Null

# Force a run boundary around the freshly inserted text (without changing
# its final look) so it stays its own <w:r>, distinct from the run that
# used to end in ";" and from the run that will hold the new ";".
$midRange2 = $d.Range($insertPos, $insertPos + $newText.Length)
$midRange2.Font.Bold = 1
$midRange2.Font.Bold = 0

# Likewise, give the trailing ";" its own run by nudging its formatting
# the same way, leaving the visible text/formatting unchanged.
$semiPos = $insertPos + $newText.Length
$semiRange = $d.Range($semiPos, $semiPos + 1)
$semiRange.Font.Bold = 1
$semiRange.Font.Bold = 0
